# covid-manual-excel.xlsx — "Update 2020-03-16: including US again."
$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy\-mm\-dd\ hh:mm"
$newDate = 43906.3375

# Rename "United States" -> "US" (keeps sheetId/rId intact)
$wsUS = $wb.Worksheets.Item("United States")
$wsUS.Name = "US"

# --- Italy (sheet1.xml) ---
$ws = $wb.Worksheets.Item("Italy")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 24747
$ws.Range("C5").Value = 2335
$ws.Range("D5").Value = 1809

# --- Germany (sheet2.xml) ---
$ws = $wb.Worksheets.Item("Germany")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 5813
$ws.Range("C5").Value = 46
$ws.Range("D5").Value = 13

# --- France (sheet3.xml) ---
$ws = $wb.Worksheets.Item("France")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 5423
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 127

# --- Spain (sheet4.xml) ---
$ws = $wb.Worksheets.Item("Spain")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 7844
$ws.Range("C5").Value = 517
$ws.Range("D5").Value = 292

# --- United Kingdom (sheet5.xml) ---
$ws = $wb.Worksheets.Item("United Kingdom")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 1372
$ws.Range("C5").Value = 18
$ws.Range("D5").Value = 35
$ws.Range("A6").NumberFormat = $dateFmt

# --- US, formerly United States (sheet6.xml) ---
$ws = $wb.Worksheets.Item("US")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 3791
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 65

# --- Austria (sheet7.xml) ---
$ws = $wb.Worksheets.Item("Austria")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 860
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 1

# --- Switzerland (sheet8.xml) ---
$ws = $wb.Worksheets.Item("Switzerland")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 2200
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 14

# --- Korea, South (sheet9.xml) ---
$ws = $wb.Worksheets.Item("Korea, South")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 8236
$ws.Range("C5").Value = 1137
$ws.Range("D5").Value = 75

# --- Iran (sheet10.xml) ---
$ws = $wb.Worksheets.Item("Iran")
$ws.Range("A5").Value = $newDate
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = 13938
$ws.Range("C5").Value = 4590
$ws.Range("D5").Value = 724

# --- Now set per-sheet selections. The LAST sheet selected/activated here
#     becomes the active tab in the workbook, so US must be last to match
#     the target (tabSelected on US's sheet, activeTab index 5). ---
$wb.Worksheets.Item("Italy").Range("D6").Select()
$wb.Worksheets.Item("Germany").Range("E5").Select()
$wb.Worksheets.Item("France").Range("E5").Select()
$wb.Worksheets.Item("Spain").Range("E5").Select()
$wb.Worksheets.Item("United Kingdom").Range("E5").Select()
$wb.Worksheets.Item("Austria").Range("E5").Select()
$wb.Worksheets.Item("Switzerland").Range("E5").Select()
$wb.Worksheets.Item("Korea, South").Range("E5").Select()
$wb.Worksheets.Item("Iran").Range("E5").Select()
$wb.Worksheets.Item("US").Range("E5").Select()

Write-Host "done"
